$wb = $excel.ActiveWorkbook

# ---- contact_data ----
$ws = $wb.Worksheets.Item('contact_data')
$ws.Range('A4').Value = 1
$ws.Range('B4').Value = '2024-12-02 15:09:35'
$ws.Range('E4').Value = 'SKZ'
$ws.Range('A5').Value = 1
$ws.Range('B5').Value = '2024-12-03 07:51:30'
$ws.Range('E5').Value = 'SKZ'
$ws.Range('A6').Value = 1
$ws.Range('B6').Value = '2024-12-03 07:55:02'
$ws.Range('E6').Value = 'SKZ'
$ws.Range('A7').Value = 1
$ws.Range('B7').Value = '2024-12-03 07:55:05'
$ws.Range('E7').Value = 'SKZ'
$ws.Range('A8').Value = 1
$ws.Range('B8').Value = '2024-12-03 07:55:08'
$ws.Range('E8').Value = 'SKZ'
$ws.Range('A9').Value = 1
$ws.Range('B9').Value = '2024-12-03 07:55:08'
$ws.Range('E9').Value = 'SKZ'
$ws.Range('A10').Value = 1
$ws.Range('B10').Value = '2024-12-03 12:56:33'
$ws.Range('E10').Value = 'SKZ'
$ws.Range('A11').Value = 1
$ws.Range('B11').Value = '2024-12-03 12:56:34'
$ws.Range('E11').Value = 'SKZ'

# ---- company_data ----
$ws = $wb.Worksheets.Item('company_data')
$ws.Range('A4').Value = 1
$ws.Range('B4').Value = '2024-12-02 15:10:07'
$ws.Range('C4').Value = 'SKZ'
$ws.Range('D4').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E4').Value = '''97076'
$ws.Range('F4').Value = 'Würzburg'
$ws.Range('G4').Value = 'Bayern'
$ws.Range('H4').Value = 'Deutschland'
$ws.Range('I4').Value = $false
$ws.Range('J4').Value = $false
$ws.Range('L4').Value = 49.80282025
$ws.Range('M4').Value = 10.00010726291456
$ws.Range('A5').Value = 1
$ws.Range('B5').Value = '2024-12-03 07:51:55'
$ws.Range('C5').Value = 'SKZ'
$ws.Range('D5').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E5').Value = '''97076'
$ws.Range('F5').Value = 'Würzburg'
$ws.Range('G5').Value = 'Bayern'
$ws.Range('H5').Value = 'Deutschland'
$ws.Range('I5').Value = $true
$ws.Range('J5').Value = $false
$ws.Range('L5').Value = 49.80282025
$ws.Range('M5').Value = 10.00010726291456
$ws.Range('A6').Value = 1
$ws.Range('B6').Value = '2024-12-03 07:55:38'
$ws.Range('C6').Value = 'SKZ'
$ws.Range('D6').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E6').Value = '''97076'
$ws.Range('F6').Value = 'Würzburg'
$ws.Range('G6').Value = 'Bayern'
$ws.Range('H6').Value = 'Deutschland'
$ws.Range('I6').Value = $true
$ws.Range('J6').Value = $false
$ws.Range('L6').Value = 49.80282025
$ws.Range('M6').Value = 10.00010726291456
$ws.Range('A7').Value = 1
$ws.Range('B7').Value = '2024-12-03 08:08:41'
$ws.Range('D7').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E7').Value = '''97076'
$ws.Range('F7').Value = 'Würzburg'
$ws.Range('G7').Value = 'Bayern'
$ws.Range('H7').Value = 'Deutschland'
$ws.Range('I7').Value = $false
$ws.Range('J7').Value = $false
$ws.Range('L7').Value = 49.80282025
$ws.Range('M7').Value = 10.00010726291456
$ws.Range('A8').Value = 1
$ws.Range('B8').Value = '2024-12-03 10:00:51'
$ws.Range('D8').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E8').Value = '''97076'
$ws.Range('F8').Value = 'Würzburg'
$ws.Range('G8').Value = 'Bayern'
$ws.Range('H8').Value = 'Deutschland'
$ws.Range('I8').Value = $false
$ws.Range('J8').Value = $false
$ws.Range('L8').Value = 49.80282025
$ws.Range('M8').Value = 10.00010726291456
$ws.Range('A9').Value = 1
$ws.Range('B9').Value = '2024-12-03 10:01:58'
$ws.Range('D9').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E9').Value = '''97076'
$ws.Range('F9').Value = 'Würzburg'
$ws.Range('G9').Value = 'Bayern'
$ws.Range('H9').Value = 'Deutschland'
$ws.Range('I9').Value = $false
$ws.Range('J9').Value = $false
$ws.Range('L9').Value = 49.80282025
$ws.Range('M9').Value = 10.00010726291456
$ws.Range('A10').Value = 1
$ws.Range('B10').Value = '2024-12-03 10:04:08'
$ws.Range('D10').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E10').Value = '''97076'
$ws.Range('F10').Value = 'Würzburg'
$ws.Range('G10').Value = 'Bayern'
$ws.Range('H10').Value = 'Deutschland'
$ws.Range('I10').Value = $false
$ws.Range('J10').Value = $false
$ws.Range('L10').Value = 49.80282025
$ws.Range('M10').Value = 10.00010726291456
$ws.Range('A11').Value = 1
$ws.Range('B11').Value = '2024-12-03 12:56:57'
$ws.Range('C11').Value = 'SKZ'
$ws.Range('D11').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E11').Value = '''97076'
$ws.Range('F11').Value = 'Würzburg'
$ws.Range('G11').Value = 'Bayern'
$ws.Range('H11').Value = 'Deutschland'
$ws.Range('I11').Value = $true
$ws.Range('J11').Value = $false
$ws.Range('L11').Value = 49.80282025
$ws.Range('M11').Value = 10.00010726291456
$ws.Range('A12').Value = 1
$ws.Range('B12').Value = '2024-12-03 12:59:06'
$ws.Range('C12').Value = 'SKZ'
$ws.Range('D12').Value = 'Friedrich-Bergius-Ring 22'
$ws.Range('E12').Value = '''97076'
$ws.Range('F12').Value = 'Würzburg'
$ws.Range('G12').Value = 'Bayern'
$ws.Range('H12').Value = 'Deutschland'
$ws.Range('I12').Value = $true
$ws.Range('J12').Value = $false
$ws.Range('L12').Value = 49.80282025
$ws.Range('M12').Value = 10.00010726291456

# ---- product_fractions ----
$ws = $wb.Worksheets.Item('product_fractions')
$ws.Range('A6').Value = 1
$ws.Range('B6').Value = '2024-12-02 15:10:27'
$ws.Range('C6').Value = '[''PE-LD'', ''PVC-P'', ''Kupfer'']'
$ws.Range('D6').Value = '['''', '''', '''']'
$ws.Range('E6').Value = '[20.0, 50.0, 30.0]'
$ws.Range('A7').Value = 1
$ws.Range('B7').Value = '2024-12-03 07:52:20'
$ws.Range('C7').Value = '[''PS'', ''PE-HD'', ''Chrom'']'
$ws.Range('D7').Value = '['''', '''', '''']'
$ws.Range('E7').Value = '[80.0, 10.0, 10.0]'
$ws.Range('A8').Value = 1
$ws.Range('B8').Value = '2024-12-03 07:56:02'
$ws.Range('C8').Value = '[''PE-LD'', ''PEEK'', ''Chrom'']'
$ws.Range('D8').Value = '['''', '''', '''']'
$ws.Range('E8').Value = '[50.0, 30.0, 20.0]'
$ws.Range('A9').Value = 1
$ws.Range('B9').Value = '2024-12-03 08:09:00'
$ws.Range('C9').Value = '[''PP'', ''Aluminium'', ''PS'']'
$ws.Range('D9').Value = '['''', '''', '''']'
$ws.Range('E9').Value = '[80.0, 10.0, 10.0]'
$ws.Range('A10').Value = 1
$ws.Range('B10').Value = '2024-12-03 10:07:06'
$ws.Range('C10').Value = '[''PE-LD'', ''PP'', ''Aluminium'']'
$ws.Range('D10').Value = '['''', '''', '''']'
$ws.Range('E10').Value = '[80.0, 10.0, 10.0]'
$ws.Range('A11').Value = 1
$ws.Range('B11').Value = '2024-12-03 12:59:39'
$ws.Range('C11').Value = '[''PE-LD'', ''ABS'', ''Chrom'']'
$ws.Range('D11').Value = '['''', '''', '''']'
$ws.Range('E11').Value = '[80.0, 10.0, 10.0]'

# ---- product_origin ----
$ws = $wb.Worksheets.Item('product_origin')
$ws.Range('A4').Value = 1
$ws.Range('B4').Value = '2024-12-02 15:10:38'
$ws.Range('C4').Value = 'Post-Industrial (PI)'
$ws.Range('D4').Value = 'Spritzguss'
$ws.Range('A5').Value = 1
$ws.Range('B5').Value = '2024-12-03 07:52:29'
$ws.Range('C5').Value = 'Post-Industrial (PI)'
$ws.Range('D5').Value = 'Spritzguss'
$ws.Range('A6').Value = 1
$ws.Range('B6').Value = '2024-12-03 10:07:15'
$ws.Range('C6').Value = 'Post-Industrial (PI)'
$ws.Range('A7').Value = 1
$ws.Range('B7').Value = '2024-12-03 12:59:50'
$ws.Range('C7').Value = 'Post-Industrial (PI)'

# ---- product_quality ----
$ws = $wb.Worksheets.Item('product_quality')
$ws.Range('A6').Value = 1
$ws.Range('B6').Value = '2024-12-02 15:10:57'
$ws.Range('C6').Value = 'Ja'
$ws.Range('D6').Value = 'bunt'
$ws.Range('E6').Value = 100
$ws.Range('F6').Value = 'gering'
$ws.Range('H6').Value = '[[''Gleitmittel''], [''Flammschutzmittel''], [''Biozide'']]'
$ws.Range('I6').Value = '[[], [], []]'
$ws.Range('A7').Value = 1
$ws.Range('B7').Value = '2024-12-03 07:56:19'
$ws.Range('C7').Value = 'Ja'
$ws.Range('D7').Value = 'grün'
$ws.Range('E7').Value = 100
$ws.Range('F7').Value = 'gering'
$ws.Range('H7').Value = '[[], [], []]'
$ws.Range('I7').Value = '[[], [], []]'
$ws.Range('A8').Value = 1
$ws.Range('B8').Value = '2024-12-03 10:07:33'
$ws.Range('C8').Value = 'Ja'
$ws.Range('D8').Value = 'braun'
$ws.Range('E8').Value = 100
$ws.Range('H8').Value = '[[], [], []]'
$ws.Range('I8').Value = '[[], [], []]'
$ws.Range('A9').Value = 1
$ws.Range('B9').Value = '2024-12-03 13:01:40'
$ws.Range('C9').Value = 'Ja'
$ws.Range('D9').Value = 'bunt'
$ws.Range('E9').Value = 100
$ws.Range('F9').Value = 'gering'
$ws.Range('H9').Value = '[[''Gleitmittel''], [], []]'
$ws.Range('I9').Value = '[[], [''Glimmer''], []]'

# ---- additive_quality ----
$ws = $wb.Worksheets.Item('additive_quality')
$ws.Range('A4').Value = 1
$ws.Range('B4').Value = '2024-12-02 15:11:40'
$ws.Range('C4').Value = '[[''Gleitmittel''], [''Flammschutzmittel''], [''Biozide'']]'
$ws.Range('D4').Value = '[[1.0], [1.0], [1.0]]'
$ws.Range('E4').Value = '[[], [], []]'
$ws.Range('F4').Value = '[[], [], []]'
$ws.Range('A5').Value = 1
$ws.Range('B5').Value = '2024-12-03 13:01:48'
$ws.Range('C5').Value = '[[''Gleitmittel''], [], []]'
$ws.Range('D5').Value = '[[5.0], [5.0], [5.0]]'
$ws.Range('E5').Value = '[[], [''Glimmer''], []]'
$ws.Range('F5').Value = '[[], [1.0], [1.0]]'

# ---- product_amount ----
$ws = $wb.Worksheets.Item('product_amount')
$ws.Range('A7').Value = 1
$ws.Range('B7').Value = '2024-12-02 15:11:48'
$ws.Range('C7').Value = 7
$ws.Range('D7').Value = 7
$ws.Range('E7').Value = 'Woche'
$ws.Range('A8').Value = 1
$ws.Range('B8').Value = '2024-12-03 07:52:56'
$ws.Range('C8').Value = 5
$ws.Range('D8').Value = 4
$ws.Range('E8').Value = 'Quartal'
$ws.Range('A9').Value = 1
$ws.Range('B9').Value = '2024-12-03 07:56:32'
$ws.Range('C9').Value = 5
$ws.Range('D9').Value = 6
$ws.Range('E9').Value = 'Monat'
$ws.Range('A10').Value = 1
$ws.Range('B10').Value = '2024-12-03 08:09:15'
$ws.Range('C10').Value = 5
$ws.Range('D10').Value = 5
$ws.Range('E10').Value = 'Monat'
$ws.Range('A11').Value = 1
$ws.Range('B11').Value = '2024-12-03 10:07:43'
$ws.Range('C11').Value = 5
$ws.Range('D11').Value = 10
$ws.Range('E11').Value = 'Quartal'
$ws.Range('A12').Value = 1
$ws.Range('B12').Value = '2024-12-03 13:02:02'
$ws.Range('C12').Value = 4.99
$ws.Range('D12').Value = 5
$ws.Range('E12').Value = 'Woche'
$ws.Range('F6').Copy($ws.Range('F12'))
$ws.Range('F6').ClearContents()

